$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (001294033 / VIVIANE / 68960.97) is removed entirely, shifting
# everything below it up by one row.
$ws.Rows.Item(2).Delete()

# The row that used to be row 5 (008115927 / ARI / 28000) is now row 4
# after the deletion above; update its values to the new account.
# A leading apostrophe forces the numeric-looking account number to be
# stored as text (same as typing it directly in Excel), preserving the
# leading zeros instead of it being auto-converted to a number.
$ws.Cells.Item(4, 1).Value = "'005624730"
$ws.Cells.Item(4, 2).Value = "ISABEL"
$ws.Cells.Item(4, 3).Value = 30000
